# Swap the presentation's theme palette from the "Integral" scheme to the
# default "Office Theme" scheme (ppt/theme/theme1.xml gains the colours that
# used to live in ppt/theme/theme2.xml, and vice versa in the source OOXML).
#
# The PowerPoint object model reaches the deck's theme colours through
# SlideMaster.ColorScheme (12 slots: dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) -- exactly the values that differ between the two theme parts;
# the font scheme and format scheme are identical between them, so driving
# ColorScheme is sufficient to reproduce the colour swap.

function Convert-HexToComRGB {
    param([string]$Hex)
    $r = [Convert]::ToInt32($Hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target values, in ColorScheme.Item(1..12) order, taken from the "Office
# Theme" clrScheme that the diff moves into theme1.xml.
$officeThemeHex = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$p = $ppt.ActivePresentation

$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

for ($i = 0; $i -lt $officeThemeHex.Count; $i++) {
    $colorScheme.Item($i + 1).RGB = Convert-HexToComRGB $officeThemeHex[$i]
}
